$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.569.87'
$ws.Range("E2").Value = '  -0.31%  '
$ws.Range("D3").Value = '1.841.24'
$ws.Range("E3").Value = '  -0.35%  '
$ws.Range("D4").Value = '1.009'
$ws.Range("E4").Value = '  -2.15%  '
$ws.Range("D5").Value = '317.09'
$ws.Range("E5").Value = '  -1.43%  '
$ws.Range("D6").Value = '1.007'
$ws.Range("E6").Value = '  -2.25%  '
$ws.Range("E7").Value = '  -1.88%  '
$ws.Range("D8").Value = '0.3726'
$ws.Range("E8").Value = '  -1.71%  '
$ws.Range("D9").Value = '0.07285'
$ws.Range("E9").Value = '  -1.27%  '
$ws.Range("D10").Value = '0.8702'
$ws.Range("E10").Value = '  -1.33%  '
$ws.Range("D11").Value = '21.32'
$ws.Range("E11").Value = '  -0.94%  '
$ws.Range("D12").Value = '1.843.10'
$ws.Range("E12").Value = '  -0.79%  '
$ws.Range("D13").Value = '6.712'
$ws.Range("D14").Value = '5.388'
$ws.Range("E14").Value = '  -1.93%  '
$ws.Range("D15").Value = '0.07111'
$ws.Range("E15").Value = '  -0.51%  '
$ws.Range("D16").Value = '88.58'
$ws.Range("E16").Value = '  +4.37%  '
$ws.Range("D17").Value = '1.010'
$ws.Range("E17").Value = '  -2.56%  '
$ws.Range("D18").Value = '0.000008965'
$ws.Range("E18").Value = '  -0.87%  '
$ws.Range("E19").Value = '  -2.24%  '
$ws.Range("D20").Value = '15.32'
$ws.Range("E20").Value = '  -0.82%  '
$ws.Range("D21").Value = '27.580.44'
$ws.Range("E21").Value = '  -0.35%  '
$ws.Range("D22").Value = '5.186'
$ws.Range("E22").Value = '  -2.02%  '
$ws.Range("E23").Value = '  -2.62%  '
$ws.Range("D24").Value = '2.071.95'
$ws.Range("E25").Value = '  -4.13%  '
$ws.Range("D26").Value = '154.38'
$ws.Range("E26").Value = '  -2.61%  '
$ws.Range("D27").Value = '18.53'
$ws.Range("E27").Value = '  -0.76%  '
$ws.Range("E28").Value = '  +8.06%  '
$ws.Range("D29").Value = '5.306'
$ws.Range("E29").Value = '  -0.10%  '
$ws.Range("D30").Value = '117.50'
$ws.Range("E30").Value = '  -0.11%  '
$ws.Range("D31").Value = '0.08896'
$ws.Range("E31").Value = '  -1.73%  '
$ws.Range("D32").Value = '1.212'
$ws.Range("E32").Value = '  +0.61%  '
$ws.Range("D33").Value = '0.7717'
$ws.Range("E33").Value = '  +0.19%  '
$ws.Range("D34").Value = '4.510'
$ws.Range("E34").Value = '  -0.89%  '
$ws.Range("D35").Value = '2.908'
$ws.Range("E35").Value = '  -3.22%  '
$ws.Range("E36").Value = '  -2.39%  '
$ws.Range("D37").Value = '1.127'
$ws.Range("E37").Value = '  -1.76%  '
$ws.Range("D38").Value = '0.01966'
$ws.Range("E38").Value = '  -0.26%  '
$ws.Range("D39").Value = '0.05291'
$ws.Range("E39").Value = '  +0.48%  '
$ws.Range("D40").Value = '2.885'
$ws.Range("E40").Value = '  +1.81%  '
$ws.Range("D41").Value = '7.131'
$ws.Range("E41").Value = '  +3.90%  '
$ws.Range("D42").Value = '0.1684'
$ws.Range("E42").Value = '  +0.99%  '
$ws.Range("D43").Value = '0.5108'
$ws.Range("E43").Value = '  -1.33%  '
$ws.Range("D44").Value = '8.744'
$ws.Range("E44").Value = '  +0.37%  '
$ws.Range("D45").Value = '10.67'
$ws.Range("E45").Value = '  -0.25%  '
$ws.Range("D46").Value = '106.87'
$ws.Range("E46").Value = '  -3.08%  '
$ws.Range("D47").Value = '0.4737'
$ws.Range("E47").Value = '  +1.04%  '
$ws.Range("D48").Value = '0.06451'
$ws.Range("E48").Value = '  -2.27%  '
$ws.Range("E49").Value = '  -2.48%  '
$ws.Range("D50").Value = '1.677'
$ws.Range("E50").Value = '  -1.27%  '
$ws.Range("D51").Value = '1.838'
$ws.Range("E51").Value = '  -3.00%  '
